# stagePlan.docx fix: bulk upload
#  1. Nudge three table column widths by 1 twip (1718->1717, 2960->2961, 2605->2604)
#  2. Merge "اجمالي المختبرين" + trailing-space run into a single run
#  3. Remove the anchored logo picture from the page header
#  4. Merge "خطة الاستخدام اليومي للمرحلة" + trailing-space run into a single run

$d = $word.ActiveDocument

# --- 1. Table column width tweaks -----------------------------------------
$t = $d.Tables(1)
$t.Columns(2).Width = 85.85   # 1717 dxa (was 1718)
$t.Columns(5).Width = 148.05  # 2961 dxa (was 2960)
$t.Columns(7).Width = 130.2   # 2604 dxa (was 2605)

# --- 2. Merge the "اجمالي المختبرين" + " " runs in the body table ---------
$d.Content.Find.Execute("اجمالي المختبرين ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "اجمالي المختبرين ", 2) | Out-Null

# --- 3. Drop the anchored logo image from the header -----------------------
$hdr = $d.Sections(1).Headers(1)
for ($i = $hdr.Shapes.Count; $i -ge 1; $i--) {
    $hdr.Shapes($i).Delete()
}

# --- 4. Merge the "خطة الاستخدام اليومي للمرحلة" + " " runs in the header -
$hdr.Range.Find.Execute("خطة الاستخدام اليومي للمرحلة ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "خطة الاستخدام اليومي للمرحلة ", 2) | Out-Null
